$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) affected by reordering of the country list ---
$ws.Range("A1").Value = 'Datos actualizados a 28 de Marzo de 2020 a las 06:29'
$ws.Range("A36").Value = 'Tailandia'
$ws.Range("A37").Value = 'Sudafrica'
$ws.Range("A43").Value = 'India'
$ws.Range("A44").Value = 'Islandia'
$ws.Range("A102").Value = 'Uzbekistan'
$ws.Range("A103").Value = 'Costa de Marfil'
$ws.Range("A104").Value = 'Camboya'
$ws.Range("A105").Value = 'Honduras'
$ws.Range("A106").Value = 'Mauricio'
$ws.Range("A107").Value = 'Bielorrusia'
$ws.Range("A108").Value = 'Martinica'
$ws.Range("A109").Value = 'Camerun'
$ws.Range("A110").Value = 'Estado de Palestina'

# --- Update statistic values (columns B-H) for rows with new data ---
$ws.Range("B4").Value = 104256
$ws.Range("C4").Value = 130
$ws.Range("E4").Value = 100027
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 1704

$ws.Range("B36").Value = 1245
$ws.Range("C36").Value = 109
$ws.Range("D36").Value = 97
$ws.Range("E36").Value = 1142
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 6

$ws.Range("B37").Value = 1170
$ws.Range("D37").Value = 31
$ws.Range("E37").Value = 1138
$ws.Range("F37").Value = 7
$ws.Range("H37").Value = 1

$ws.Range("B43").Value = 902
$ws.Range("C43").Value = 15
$ws.Range("D43").Value = 83
$ws.Range("E43").Value = 799
$ws.Range("F43").Value = 0
$ws.Range("H43").Value = 20

$ws.Range("B44").Value = 890
$ws.Range("D44").Value = 97
$ws.Range("E44").Value = 791
$ws.Range("F44").Value = 18
$ws.Range("H44").Value = 2

$ws.Range("B102").Value = 104
$ws.Range("C102").Value = 16
$ws.Range("D102").Value = 5
$ws.Range("E102").Value = 97
$ws.Range("F102").Value = 8
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 2

$ws.Range("B103").Value = 101
$ws.Range("D103").Value = 3
$ws.Range("E103").Value = 98
$ws.Range("F103").Value = 0

$ws.Range("B104").Value = 99
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 13
$ws.Range("E104").Value = 86
$ws.Range("F104").Value = 1
$ws.Range("H104").Value = 0

$ws.Range("B105").Value = 95
$ws.Range("C105").Value = 27
$ws.Range("D105").Value = 3
$ws.Range("E105").Value = 91
$ws.Range("F105").Value = 4
$ws.Range("H105").Value = 1

$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 92
$ws.Range("F106").Value = 1

$ws.Range("B107").Value = 94
$ws.Range("D107").Value = 32
$ws.Range("E107").Value = 62
$ws.Range("F107").Value = 2
$ws.Range("H107").Value = 0

$ws.Range("B108").Value = 93
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 92
$ws.Range("F108").Value = 12
$ws.Range("H108").Value = 1

$ws.Range("D109").Value = 2
$ws.Range("E109").Value = 87
$ws.Range("H109").Value = 2

$ws.Range("B110").Value = 91
$ws.Range("D110").Value = 17
$ws.Range("E110").Value = 73
$ws.Range("F110").Value = 0

